# Auto-generated edit script: splits multi-sentiment freetext rows into one row per sentiment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(18, 10276361, 'Kevin at KEF desk was very helpful and helped me to not feel rushed through the process of check in. ', 'positive')
    ,@(18, 10276361, 'Personal item storage under the seat in front was very low to the floor compared to other planes.', 'negative')
    ,@(42, 10284760, 'The Passport Control process on Jan 1 was poorly organized. I realize that the airport controls much of this, and the wind conditions made all the planes'' de boarding occur at once.  Still the hour + wait (there as a fistfight further in front of us), followed by a 30 min. "Random security check" left me overheated and anxious on the 5 hr flight to Boston; plus, you were out of food by the time the cart got back to us. ', 'negative')
    ,@(42, 10284760, 'In contrast, I successfully upgraded to Saga class for the previous leg from CPH to KEF, which got us into the SAS Lounge for the unexpected 6 hr layover at CPH, after you cancelled service from Helsinki.  Service and food in Saga on the CPH flight were great; seats reminded me of 1992 (seemed to be the age of the plane…)', 'positive')
    ,@(81, 10268463, 'Baggage drop off took very long, could be improved. ', 'negative')
    ,@(81, 10268463, 'Fantastic: a member of staff passed the queue and helped preparing covid procedures in advance', 'positive')
    ,@(115, 10274655, 'The staff on board were excellent. ', 'positive')
    ,@(115, 10274655, 'Check in was very good at Manchester but exceptional at KEF.', 'negative')
    ,@(134, 10287707, 'Queuing at the gate in Iceland Airport with no seating was inconvenient', 'negative')
    ,@(134, 10287707, ' However one of the best flight experiences I have had for along time', 'positive')
    ,@(177, 10267000, 'Smoothest landing I''ve ever experienced.', 'positive')
    ,@(177, 10267000, 'You should have served lunch with the price of the ticket.', 'negative')
    ,@(220, 10279808, '- The headrests slide down after being placed up. ', 'negative')
    ,@(220, 10279808, '- The new 2021 in-flight entertainment has better movies than the 2020 options.', 'positive')
    ,@(220, 10279808, '- The baggage check-in at JFK airport took a long time (almost an hour).', 'negative')
    ,@(264, 10294976, 'Just to repeat all I said about the airport check in staff.   Very kind, helpful and courteous.  We flew from Reykjavik to Heathrow on 3 January 7:40am', 'positive')
    ,@(264, 10294976, 'One thing was a little annoying which was waiting for quite a long time for the bus to the plane in a narrow corridor', 'negative')
    ,@(402, 10295336, 'One flight was cancelled but the other flight was such a good experience it overrode the cancelled flight', 'negative')
    ,@(402, 10295336, ' Just waiting for my refund at this point', 'positive')
    ,@(411, 10297384, 'Please turn volume up/louder when flight attendant make announcements.  It was difficult to hear what was being said. Some of your client might have a hearing disability.', 'negative')
    ,@(411, 10297384, 'Noticed some announcement was displayed on the screens/monitors. That was helpful.  But is it possible to display all announcement of screens?', 'positive')
    ,@(413, 10297006, 'Wifi onboard was really bad', 'negative')
    ,@(413, 10297006, ' Other then that it was a great flight', 'positive')
    ,@(431, 10297711, '3of our 4 flights were delayed, which made the transfers/ connection stressful. ', 'negative')
    ,@(431, 10297711, 'Love the communication updates from', 'positive')
    ,@(431, 10297711, 'The pilot throughout the flight. Love the speed and efficiency of passport control', 'negative')
    ,@(459, 10298580, 'It''s much, much, much too hot on the plane', 'negative')
    ,@(459, 10298580, ' Seats are too small', 'negative')
    ,@(459, 10298580, ' Food is good', 'positive')
    ,@(459, 10298580, ' Service is excellent', 'positive')
    ,@(462, 10297844, 'Not sure if this was due to staff shortages because of covid, but opening check-in at Charles De Gaulle airport only 2.5 hours before take-off is terrible. We arrived 4 hours early but had to wait for check-in to open, then the line was so long that it took 2 hours to get to the front of the line, and we had to run to catch the plane. The plane thankfully waited for all passengers, but that delayed takeoff by long enough that we had to run to catch our connecting flight in Reykjavik. So, no time to buy water or food at either airport, and since there is very limited and expensive food service on board, we had to stay thirsty and hungry longer than was comfortable.', 'negative')
    ,@(462, 10297844, 'The crew was absolutely efficient and friendly as always, and we love flying Icelandiar otherwise. We are looking forward to fly again with you, hopefully when the pandemic subsides.', 'positive')
    ,@(466, 10299322, 'My overall flight experience was very easy. The seats are comfortable enough in economy and the in-flight entertainment was good. The layover in Iceland is great and there are plenty of food and beverage options at the airport. During the flights, the cabin crew was attentive during service, but since there isn''t a meal service (unless you''ve pre-paid for a meal or buy something onboard), it feels like there''s a long wait in between something as simple as water. I know changes have been made for everyone''s safety with Covid regarding service- but perhaps offering a water bottle in lieu of very small and occasional cups of water would be an option?', 'positive')
    ,@(466, 10299322, 'My one very minor qualm is that I did purchase the in-flight wifi and it was inconsistent for something as simple as checking emails... I certainly don''t mind paying for a service, but I would hope that I am getting what I am paying an additional fee for.', 'negative')
    ,@(484, 10306628, 'Arrival was fluency even with the management for Covid 19. Get off the plane, security, pick up the luggages , covid test all was well managed easy and fast.', 'positive')
    ,@(484, 10306628, 'The icelandair agents were all professionals with courtesy. ', 'positive')
    ,@(484, 10306628, 'The problem is when we need help with a partner''s section. Lack of assistance line 24/24 for the short time departures . The planes fly all the day and whatever the time zone.', 'negative')
    ,@(561, 10305669, 'The seats on the 767 are very uncomfortable for anyone with back problems. ', 'negative')
    ,@(561, 10305669, 'Icelandair crew have only been kind and excellent, wonderful people!', 'positive')
    ,@(642, 10326492, 'We were thankful at the seat leg room and size were not as bad as United. ', 'positive')
    ,@(642, 10326492, 'It would really be helpful with a meal on board the over alantic flight and not just a pizza. When we had no transfer time we were unable to buy food for our dietary needs in the airport. Diabetic people must eat on regular basis, so it was very disappointing not offered a proper meal.', 'negative')
    ,@(729, 10316963, 'There should at least be some basic complimentary beverages provided', 'negative')
    ,@(729, 10316963, ' Staff was very friendly and kind', 'positive')
    ,@(742, 10315718, 'Yes. Everything was great except: There was no food on a five-hour flight, not good. And on my way back, onboarding we were made to walk in the snow and wind without warning to the plane. The gate was also tiny and we were not boarded properly from back to front, so many of us were waiting in the cold outside to get on the plane.', 'negative')
    ,@(742, 10315718, 'I did like the service and I thought Icelandair was great for its short travel time.', 'positive')
    ,@(771, 10318342, 'The SEA-KEF leg of the trip was great. We flew Saga Premium and had a great experience. The check-in/baggage drop off portion was the easiest and quickest I''ve experienced. The flight attendants were lovely- I forgot to order a vegetarian meal and they gave me a variety of vegetarian options they had, everything felt really clean, and it would 10/10 recommend. ', 'positive')
    ,@(771, 10318342, 'The KEF-ARN was a very different experience. We (my partner and a friend) were in row 9 seats A, B, and C. In 9D, there was a passenger who barely wore a mask and had a lot of Covid symptoms. He was coughing/hacking quite a bit, was blowing his nose frequently, and then would pull his mask under his chin. I asked him to please wear a mask and he moved it up for awhile, but eventually pulled it down again. I went to talk to a flight attendant shortly after we departed and asked if there was any way we could pay to move to Saga Premium because the guy close to us wasn''t wearing a mask and I had already asked him, he would pull it up before the flight attendants walked by and then pull it down again, he was coughing and blowing his nose frequently, and it was just quite uncomfortable and was making us feel really claustrophobic. The flight attendant said the flight was pretty full so we couldn''t change and we also couldn''t move/pay to change to the nearly empty first class seats, but she would talk to him. She did talk to him, I believe she gave him a new mask and tissues, and he wore his mask again for awhile, but mostly took it off. The flight attendants walked by him several times and didn''t say anything (which to some extent I understand because of how many viral videos these days there are of passengers acting terribly and getting into fights), but it was just discouraging. The plane was quite warm and with the mostly maskless sick passenger close to us it just made for a rather unpleasant flight.', 'negative')
    ,@(781, 10318727, 'The airline experience', 'positive')
    ,@(781, 10318727, 'Was great but I was supervised that the airline could not compensate for the extreme inconvenience and expense I was caused.  God Bless I know we are in tough times but more could have been done.', 'negative')
    ,@(829, 10325169, 'I wish check-in staff were more approachable and the staff who check your visa/entry requirements weren''t so abrupt with their questionning', 'negative')
    ,@(829, 10325169, ' It really left a sour taste after incredible on board service and meal, and a wonderful lounge', 'positive')
    ,@(900, 10326826, 'Reykjavik airport saga lounge''s location is strange.  It''s past the passport control and not every lounge-eligible customer can access it.  It should have access to travelers without going through passport control.  I was saga business transit passenger but couldn''t use it due this reason.  ', 'negative')
    ,@(900, 10326826, 'Otherwise great experience with Icelander airline.  Liked it!', 'positive')
    ,@(940, 10329335, 'Everything was EXCELLENT except one encounter with agent at check in at KEF.  There was no issue with my carry-on departing and it weighed less on my return. When asked to put it in the bin to size it I did. It fit completely in . He said it was to big I said it fits completely in , tightly-but it fit. Explained there was no issue coming and that we needed the clothes that were in there in case they lost our luggage or we were held over in JFK. He let us go and then when we got to the gate we were red flagged. He had sent a note to the gate to let us board with it. The agent at the gate asked me to put it in the sizing bin. It fit completely in and she didn''t know why he flagged it. ', 'positive')
    ,@(940, 10329335, 'It felt sneaky and did not leave a comfortable test with us.', 'negative')
    ,@(952, 10328617, 'Need better intertainment like new movies etc', 'negative')
    ,@(952, 10328617, 'pls provide free snacks for 7 hours flight not just 2 time drink, this is the first time I had to pay for snacks for 7 hours flight,but  flight attendance are very nice', 'positive')
)

$lastOldRow = 45
$firstRow = 2

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $firstRow + $i
    if ($r -gt $lastOldRow) {
        # Row did not exist before -- copy formatting (incl. style s="1" on col A) from row 2
        $ws.Cells.Item($firstRow, 1).Copy($ws.Cells.Item($r, 1))
    }
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

